$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q11)
$ws.Range("B7").Value = 0.08022968598094522
$ws.Range("C7").Value = 0.9401611245884408
$ws.Range("D7").Value = 1.882473483720003
$ws.Range("E7").Value = 1.372032610297584
$ws.Range("F7").Value = 1.38807073857714
$ws.Range("G7").Value = 38

# Row 8 (Q12)
$ws.Range("B8").Value = 0.09684414037220078
$ws.Range("C8").Value = 1.117279359203689
$ws.Range("D8").Value = 3.04125129278947
$ws.Range("E8").Value = 1.743918373316099
$ws.Range("F8").Value = 1.76524536140641
$ws.Range("G8").Value = 37

# Row 9 (Q13)
$ws.Range("B9").Value = 0.1091771171770269
$ws.Range("C9").Value = 1.083714881619536
$ws.Range("D9").Value = 1.693673186052097
$ws.Range("E9").Value = 1.301411997044786
$ws.Range("F9").Value = 1.330513768253694
$ws.Range("G9").Value = 20

# Row 10 (Q14)
$ws.Range("B10").Value = 0.3732545068521868
$ws.Range("C10").Value = 1.052204277644287
$ws.Range("D10").Value = 1.860053869313877
$ws.Range("E10").Value = 1.363837919004262
$ws.Range("F10").Value = 1.365331286158824
$ws.Range("G10").Value = 13

# Row 11 (Q15)
$ws.Range("B11").Value = 0.9706231273330275
$ws.Range("C11").Value = 1.152940683320901
$ws.Range("D11").Value = 1.587806877227315
$ws.Range("E11").Value = 1.260082091463614
$ws.Range("F11").Value = 0.8983997035796265
$ws.Range("G11").Value = 5
